$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G4").Value = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G7").Value = "Kerelos.zareef@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
$ws.Range("G15").Value = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G17").Value = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"
$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G30").Value = "wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
